$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate match data among rows 235, 237, 239 (A=233,235,237) — B/F/G/H/I/J/K..AC
$ws.Range("B235").Value = 6861095
$ws.Range("F235").Value = "FC Botosani"
$ws.Range("G235").Value = "Farul Constanta"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("K235").Value = 3.75
$ws.Range("L235").Value = 3.4
$ws.Range("M235").Value = 1.909
$ws.Range("N235").Value = 3.1
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.8
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = 2
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 0.3875
$ws.Range("AA235").Value = -0.5
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 1.05

$ws.Range("B237").Value = 6836277
$ws.Range("F237").Value = "CFR Cluj"
$ws.Range("G237").Value = "AFC Hermannstadt"
$ws.Range("H237").Value = 1
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = "H"
$ws.Range("K237").Value = 1.7
$ws.Range("L237").Value = 3.4
$ws.Range("M237").Value = 5
$ws.Range("N237").Value = 1.65
$ws.Range("O237").Value = 3.5
$ws.Range("P237").Value = 5.25
$ws.Range("Q237").Value = -0.75
$ws.Range("R237").Value = 1.85
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.875
$ws.Range("V237").Value = 1.975
$ws.Range("W237").Value = 0.6499999999999999
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.425
$ws.Range("AA237").Value = -0.5
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.9750000000000001

$ws.Range("B239").Value = 6852370
$ws.Range("F239").Value = "Dinamo Bucharest"
$ws.Range("G239").Value = "ACS UTA Batrana Doamna"
$ws.Range("H239").Value = 1
$ws.Range("I239").Value = 0
$ws.Range("J239").Value = "H"
$ws.Range("K239").Value = 2.55
$ws.Range("L239").Value = 2.875
$ws.Range("M239").Value = 3
$ws.Range("N239").Value = 2.375
$ws.Range("O239").Value = 3
$ws.Range("P239").Value = 3.1
$ws.Range("Q239").Value = -0.25
$ws.Range("R239").Value = 2
$ws.Range("S239").Value = 1.85
$ws.Range("T239").Value = 2.25
$ws.Range("U239").Value = 1.975
$ws.Range("V239").Value = 1.875
$ws.Range("W239").Value = 1.375
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = 1
$ws.Range("AA239").Value = -1
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.875

# Row 265 (A=263): fill in result + odds deltas
$ws.Range("H265").Value = 0
$ws.Range("I265").Value = 0
$ws.Range("J265").Value = "D"
$ws.Range("O265").Value = 3.1
$ws.Range("U265").Value = 1.9
$ws.Range("V265").Value = 1.95
$ws.Range("W265").Value = -1
$ws.Range("X265").Value = 2.1
$ws.Range("Y265").Value = -1
$ws.Range("Z265").Value = -0.5
$ws.Range("AA265").Value = 0.4125
$ws.Range("AB265").Value = -1
$ws.Range("AC265").Value = 0.95

# Row 266 (A=264): fill in result + odds deltas
$ws.Range("H266").Value = 1
$ws.Range("I266").Value = 0
$ws.Range("J266").Value = "H"
$ws.Range("O266").Value = 2.8
$ws.Range("P266").Value = 3.6
$ws.Range("T266").Value = 1.75
$ws.Range("U266").Value = 1.85
$ws.Range("V266").Value = 2
$ws.Range("W266").Value = 1.15
$ws.Range("X266").Value = -1
$ws.Range("Y266").Value = -1
$ws.Range("Z266").Value = 0.875
$ws.Range("AA266").Value = -1
$ws.Range("AB266").Value = -1
$ws.Range("AC266").Value = 1

# New upcoming fixtures: rows 267-270 (A=265..268)
$ws.Range("A266").Copy() | Out-Null
$ws.Range("A267").PasteSpecial(-4122) | Out-Null
$ws.Range("E266").Copy() | Out-Null
$ws.Range("E267").PasteSpecial(-4122) | Out-Null
$ws.Range("A267").Value = 265
$ws.Range("B267").Value = 7951792
$ws.Range("C267").Value = "Romania Liga I"
$ws.Range("D267").Value = "Romania Liga I"
$ws.Range("E267").Value = 45394.47916666666
$ws.Range("F267").Value = "ACS UTA Batrana Doamna"
$ws.Range("G267").Value = "AFC Hermannstadt"
$ws.Range("K267").Value = 2.3
$ws.Range("L267").Value = 2.9
$ws.Range("M267").Value = 3.1
$ws.Range("N267").Value = 2.5
$ws.Range("O267").Value = 2.9
$ws.Range("P267").Value = 2.875
$ws.Range("Q267").Value = 0
$ws.Range("R267").Value = 1.775
$ws.Range("S267").Value = 2.1
$ws.Range("T267").Value = 2
$ws.Range("U267").Value = 1.875
$ws.Range("V267").Value = 1.975
$ws.Range("W267").Value = 0
$ws.Range("X267").Value = 0
$ws.Range("Y267").Value = 0
$ws.Range("Z267").Value = 0
$ws.Range("AA267").Value = 0

$ws.Range("A266").Copy() | Out-Null
$ws.Range("A268").PasteSpecial(-4122) | Out-Null
$ws.Range("E266").Copy() | Out-Null
$ws.Range("E268").PasteSpecial(-4122) | Out-Null
$ws.Range("A268").Value = 266
$ws.Range("B268").Value = 7951755
$ws.Range("C268").Value = "Romania Liga I"
$ws.Range("D268").Value = "Romania Liga I"
$ws.Range("E268").Value = 45394.60416666666
$ws.Range("F268").Value = "CS U Craiova"
$ws.Range("G268").Value = "Farul Constanta"
$ws.Range("K268").Value = 1.8
$ws.Range("L268").Value = 3.4
$ws.Range("M268").Value = 4.2
$ws.Range("N268").Value = 1.8
$ws.Range("O268").Value = 3.4
$ws.Range("P268").Value = 4.2
$ws.Range("Q268").Value = -0.5
$ws.Range("R268").Value = 1.825
$ws.Range("S268").Value = 2.025
$ws.Range("T268").Value = 2.5
$ws.Range("U268").Value = 1.95
$ws.Range("V268").Value = 1.9
$ws.Range("W268").Value = 0
$ws.Range("X268").Value = 0
$ws.Range("Y268").Value = 0
$ws.Range("Z268").Value = 0
$ws.Range("AA268").Value = 0

$ws.Range("A266").Copy() | Out-Null
$ws.Range("A269").PasteSpecial(-4122) | Out-Null
$ws.Range("E266").Copy() | Out-Null
$ws.Range("E269").PasteSpecial(-4122) | Out-Null
$ws.Range("A269").Value = 267
$ws.Range("B269").Value = 7951791
$ws.Range("C269").Value = "Romania Liga I"
$ws.Range("D269").Value = "Romania Liga I"
$ws.Range("E269").Value = 45395.64583333334
$ws.Range("F269").Value = "Universitatea Cluj"
$ws.Range("G269").Value = "Petrolul Ploiesti"
$ws.Range("K269").Value = 1.909
$ws.Range("L269").Value = 3.2
$ws.Range("M269").Value = 4
$ws.Range("N269").Value = 1.909
$ws.Range("O269").Value = 3.2
$ws.Range("P269").Value = 4
$ws.Range("Q269").Value = -0.5
$ws.Range("R269").Value = 1.975
$ws.Range("S269").Value = 1.875
$ws.Range("T269").Value = 2
$ws.Range("U269").Value = 1.875
$ws.Range("V269").Value = 1.975
$ws.Range("W269").Value = 0
$ws.Range("X269").Value = 0
$ws.Range("Y269").Value = 0
$ws.Range("Z269").Value = 0
$ws.Range("AA269").Value = 0

$ws.Range("A266").Copy() | Out-Null
$ws.Range("A270").PasteSpecial(-4122) | Out-Null
$ws.Range("E266").Copy() | Out-Null
$ws.Range("E270").PasteSpecial(-4122) | Out-Null
$ws.Range("A270").Value = 268
$ws.Range("B270").Value = 7951754
$ws.Range("C270").Value = "Romania Liga I"
$ws.Range("D270").Value = "Romania Liga I"
$ws.Range("E270").Value = 45396.625
$ws.Range("F270").Value = "CFR Cluj"
$ws.Range("G270").Value = "FCSB"
$ws.Range("K270").Value = 2.3
$ws.Range("L270").Value = 3.1
$ws.Range("M270").Value = 3
$ws.Range("N270").Value = 2.3
$ws.Range("O270").Value = 3.1
$ws.Range("P270").Value = 3
$ws.Range("Q270").Value = -0.25
$ws.Range("R270").Value = 2.05
$ws.Range("S270").Value = 1.8
$ws.Range("T270").Value = 2.25
$ws.Range("U270").Value = 1.85
$ws.Range("V270").Value = 2
$ws.Range("W270").Value = 0
$ws.Range("X270").Value = 0
$ws.Range("Y270").Value = 0
$ws.Range("Z270").Value = 0
$ws.Range("AA270").Value = 0

$excel.CutCopyMode = 0

